# Move the two corporate logos on the Slide Master and refresh the
# cached "date" placeholder text (master, all layouts, handout master,
# notes master) from 03.03.2019 -> 04.03.2019.

$p = $ppt.ActivePresentation

$oldDate = "03.03.2019"
$newDate = "04.03.2019"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1) Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 2) Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# 3) Handout master date placeholder.
Update-DatePlaceholder $p.HandoutMaster.Shapes

# 4) Notes master date placeholder.
Update-DatePlaceholder $p.NotesMaster.Shapes

# 5) Reposition the two logo pictures on the slide master.
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -eq "Graphic 6") {
        $sh.Left = 26.25
        $sh.Top = 18.99992125984252
    } elseif ($sh.Name -eq "Picture 9") {
        $sh.Left = 827.9999212598425
        $sh.Top = 18.99992125984252
    }
}
